$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated odds values for rows 3, 4, 6, 7, 9 per upstream data refresh
$ws.Range("H3").Value = 3
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
$ws.Range("W3").Value = 5.5
$ws.Range("AA3").Value = 21
$ws.Range("AC3").Value = 6
$ws.Range("AG3").Value = 8.5
$ws.Range("AO3").Value = 12
$ws.Range("AU3").Value = 9.5
$ws.Range("AX3").Value = 23
$ws.Range("G4").Value = 3.25
$ws.Range("I4").Value = 2.35
$ws.Range("L4").Value = 3.25
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.7
$ws.Range("R4").Value = 1.44
$ws.Range("Z4").Value = 34
$ws.Range("AI4").Value = 11
$ws.Range("AJ4").Value = 23
$ws.Range("AP4").Value = 34
$ws.Range("AQ4").Value = 67
$ws.Range("AW4").Value = 4.33
$ws.Range("BA4").Value = 101
$ws.Range("G6").Value = 1.8
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 2.6
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7
$ws.Range("AG6").Value = 9.5
$ws.Range("AP6").Value = 29
$ws.Range("AQ6").Value = 41
$ws.Range("G7").Value = 1.38
$ws.Range("I7").Value = 8.5
$ws.Range("K7").Value = 2.3
$ws.Range("AD7").Value = 8.5
$ws.Range("AK7").Value = 51
$ws.Range("G9").Value = 2.63
$ws.Range("I9").Value = 2.35
$ws.Range("Q9").Value = 1.83
$ws.Range("R9").Value = 2.03
$ws.Range("Y9").Value = 10
$ws.Range("AB9").Value = 26
$ws.Range("AG9").Value = 9.5
$ws.Range("AH9").Value = 13
$ws.Range("BB9").Value = 126
